$wb = $excel.ActiveWorkbook

# --- Sheet "Key" (1st sheet): insert a new column A, fill with "Table"/"Data" marker ---
$wsKey = $wb.Worksheets.Item(1)
$wsKey.Columns("A:A").Insert()
$wsKey.Range("A1").Value = "Table"
$wsKey.Range("A2:A49").Value = "Data"

# --- Sheet "DWER" (2nd sheet): same treatment ---
$wsDwer = $wb.Worksheets.Item(2)
$wsDwer.Columns("A:A").Insert()
$wsDwer.Range("A1").Value = "Table"
$wsDwer.Range("A2:A41").Value = "Data"

# --- Selections: Key becomes the active/selected sheet with A2:A49 selected,
#     DWER keeps a plain A2:A41 selection (no longer the active tab) ---
[void]$wsDwer.Range("A2:A41").Select()
[void]$wsKey.Select()
[void]$wsKey.Range("A2:A49").Select()
